$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J4").Value = 2.95
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 3.45
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 1.75
$ws.Range("T4").Value = 2.75
$ws.Range("U4").Value = 1.75
$ws.Range("V4").Value = 1.98
$ws.Range("W4").Value = 8
$ws.Range("X4").Value = 12
$ws.Range("Y4").Value = 9.25
$ws.Range("AA4").Value = 19.5
$ws.Range("AB4").Value = 29
$ws.Range("AE4").Value = 13.5
$ws.Range("AG4").Value = 8.75
$ws.Range("AH4").Value = 15
$ws.Range("AI4").Value = 10.25
$ws.Range("AJ4").Value = 35
$ws.Range("AK4").Value = 25
$ws.Range("AL4").Value = 32
$ws.Range("AM4").Value = 500
$ws.Range("AP4").Value = 19
$ws.Range("AR4").Value = 80
$ws.Range("AS4").Value = 250
$ws.Range("AT4").Value = 2.75
$ws.Range("AY4").Value = 22
$ws.Range("BA4").Value = 100
$ws.Range("BB4").Value = 250
$ws.Range("H6").Value = 5.1
$ws.Range("I6").Value = 7.9
$ws.Range("K6").Value = 2.62
$ws.Range("L6").Value = 6.9
$ws.Range("S6").Value = 1.26
$ws.Range("T6").Value = 3.45
$ws.Range("U6").Value = 1.83
$ws.Range("V6").Value = 1.88
$ws.Range("AB6").Value = 24
$ws.Range("AF6").Value = 80
$ws.Range("AG6").Value = 26
$ws.Range("AH6").Value = 60
$ws.Range("AO6").Value = 5.6
$ws.Range("AT6").Value = 3.45
$ws.Range("AW6").Value = 9.25
$ws.Range("AZ6").Value = 300
$ws.Range("BB6").Value = 450
$ws.Range("G7").Value = 4.6
$ws.Range("H7").Value = 3.45
$ws.Range("I7").Value = 1.7
$ws.Range("J7").Value = 5
$ws.Range("L7").Value = 2.25
$ws.Range("N7").Value = 7
$ws.Range("Q7").Value = 1.98
$ws.Range("R7").Value = 1.78
$ws.Range("U7").Value = 1.9
$ws.Range("V7").Value = 1.8
$ws.Range("W7").Value = 11.75
$ws.Range("X7").Value = 26
$ws.Range("Y7").Value = 15.5
$ws.Range("Z7").Value = 80
$ws.Range("AA7").Value = 50
$ws.Range("AB7").Value = 55
$ws.Range("AC7").Value = 7
$ws.Range("AD7").Value = 6.9
$ws.Range("AE7").Value = 17
$ws.Range("AF7").Value = 90
$ws.Range("AG7").Value = 6.4
$ws.Range("AH7").Value = 7.5
$ws.Range("AJ7").Value = 13
$ws.Range("AK7").Value = 14
$ws.Range("AL7").Value = 29
$ws.Range("AM7").Value = 800
$ws.Range("AN7").Value = 6.4
$ws.Range("AO7").Value = 28
$ws.Range("AP7").Value = 32
$ws.Range("AQ7").Value = 175
$ws.Range("AR7").Value = 200
$ws.Range("AS7").Value = 450
$ws.Range("AU7").Value = 7.6
$ws.Range("AV7").Value = 70
$ws.Range("AW7").Value = 3.5
$ws.Range("AX7").Value = 8.25
$ws.Range("AY7").Value = 17.5
$ws.Range("AZ7").Value = 27
$ws.Range("BA7").Value = 60
$ws.Range("L8").Value = 3.45
$ws.Range("V8").Value = 2.22
$ws.Range("W8").Value = 9.5
$ws.Range("X8").Value = 13
$ws.Range("AA8").Value = 17
$ws.Range("AG8").Value = 10.25
$ws.Range("AL8").Value = 28
$ws.Range("AN8").Value = 4.4
$ws.Range("AO8").Value = 11.75
$ws.Range("AP8").Value = 17.5
$ws.Range("AX8").Value = 16
$ws.Range("BA8").Value = 100
